# Apply cryptos list update (prices / volume% / row reorders) per commit
# "Updated cryptos list on Wed Mar  1 08:46:00 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as TEXT (matches source data which stores
# prices/percentages/links as plain strings, not numbers) and then reset
# the style back to Normal so the "quote prefix" marker Excel adds for
# apostrophe-led text doesn't linger as a formatting diff.
function Set-TextValue($addr, $val) {
    $ws.Range($addr).Value = "'" + $val
    $ws.Range($addr).Style = "Normal"
}

Set-TextValue "D2" "23.749.81"
Set-TextValue "E2" "  +1.80%  "

Set-TextValue "D3" "1.653.82"
Set-TextValue "E3" "  +1.82%  "

Set-TextValue "D4" "0.9996"
Set-TextValue "E4" "  -0.20%  "

Set-TextValue "D5" "0.9997"
Set-TextValue "E5" "  -0.17%  "

Set-TextValue "D6" "304.05"
Set-TextValue "E6" "  +0.32%  "

Set-TextValue "D7" "0.3829"
Set-TextValue "E7" "  +2.39%  "

Set-TextValue "B8" "Cardano"
Set-TextValue "C8" "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
Set-TextValue "D8" "0.3614"
Set-TextValue "E8" "  -0.10%  "

Set-TextValue "B9" "OKB"
Set-TextValue "C9" "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue "D9" "51.34"
Set-TextValue "E9" "  -0.16%  "

Set-TextValue "D10" "1.252"
Set-TextValue "E10" "  +2.36%  "

Set-TextValue "D11" "0.08242"
Set-TextValue "E11" "  +1.29%  "

Set-TextValue "D12" "0.9995"
Set-TextValue "E12" "  -0.25%  "

Set-TextValue "D13" "22.77"
Set-TextValue "E13" "  +2.35%  "

Set-TextValue "D14" "6.549"
Set-TextValue "E14" "  +1.39%  "

Set-TextValue "D15" "7.418"
Set-TextValue "E15" "  +1.98%  "

Set-TextValue "D16" "0.00001236"
Set-TextValue "E16" "  -0.16%  "

Set-TextValue "D17" "1.651.70"
Set-TextValue "E17" "  +2.15%  "

Set-TextValue "D18" "97.62"
Set-TextValue "E18" "  +4.10%  "

Set-TextValue "D19" "0.06984"
Set-TextValue "E19" "  +0.58%  "

Set-TextValue "D20" "6.787"
Set-TextValue "E20" "  +3.82%  "

Set-TextValue "D21" "17.73"
Set-TextValue "E21" "  +1.33%  "

Set-TextValue "D22" "0.9996"
Set-TextValue "E22" "  -0.22%  "

Set-TextValue "D23" "12.64"
Set-TextValue "E23" "  +1.12%  "

Set-TextValue "D24" "23.739.55"
Set-TextValue "E24" "  +1.78%  "

Set-TextValue "D25" "2.536"
Set-TextValue "E25" "  +3.22%  "

Set-TextValue "D26" "3.079"
Set-TextValue "E26" "  -1.26%  "

Set-TextValue "D27" "21.34"
Set-TextValue "E27" "  +0.94%  "

Set-TextValue "D28" "151.64"
Set-TextValue "E28" "  +1.08%  "

Set-TextValue "D29" "5.269"
Set-TextValue "E29" "  +0.45%  "

Set-TextValue "D30" "135.06"
Set-TextValue "E30" "  +1.88%  "

Set-TextValue "D31" "1.835.31"
Set-TextValue "E31" "  +1.96%  "

Set-TextValue "D32" "6.883"
Set-TextValue "E32" "  +2.43%  "

Set-TextValue "D33" "1.087"
Set-TextValue "E33" "  +4.89%  "

Set-TextValue "B34" "FraxShare"
Set-TextValue "C34" "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue "D34" "11.91"
Set-TextValue "E34" "  +9.31%  "

Set-TextValue "B35" "WEMIXTOKEN"
Set-TextValue "C35" "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue "D35" "2.108"
Set-TextValue "E35" "  +0.17%  "

Set-TextValue "D36" "0.02848"
Set-TextValue "E36" "  +3.36%  "

Set-TextValue "D37" "0.2522"
Set-TextValue "E37" "  +1.17%  "

Set-TextValue "B38" "InternetComputer(DFINITY)"
Set-TextValue "C38" "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue "D38" "6.134"
Set-TextValue "E38" "  +2.83%  "

Set-TextValue "B39" "Stellar"
Set-TextValue "C39" "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue "D39" "0.08839"
Set-TextValue "E39" "  +0.76%  "

Set-TextValue "D40" "0.07058"
Set-TextValue "E40" "  -0.46%  "

Set-TextValue "D41" "12.91"
Set-TextValue "E41" "  +7.09%  "

Set-TextValue "D42" "0.7081"
Set-TextValue "E42" "  +1.64%  "

Set-TextValue "D43" "1.342"
Set-TextValue "E43" "  +0.36%  "

Set-TextValue "D44" "15.97"
Set-TextValue "E44" "  +0.31%  "

Set-TextValue "D45" "0.6563"
Set-TextValue "E45" "  +1.60%  "

Set-TextValue "E46" "  +3.57%  "

Set-TextValue "E47" "  -0.15%  "

Set-TextValue "D48" "3.974"
Set-TextValue "E48" "  +0.41%  "

Set-TextValue "D49" "0.07990"
Set-TextValue "E49" "  +0.30%  "

Set-TextValue "D50" "128.98"
Set-TextValue "E50" "  +2.77%  "

Set-TextValue "D51" "1.199"
Set-TextValue "E51" "  +1.40%  "
